$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price (D) values so Excel
# does not silently convert them to floating point numbers,
# then restore the default "Normal" style so no visual/style diff remains.

$ws.Range("D2").Value = "28.542.93"
$ws.Range("E2").Value = "  +2.21%  "

$ws.Range("D3").Value = "1.578.21"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  +0.66%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("E7").Value = "  +0.61%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.65"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.02"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.68%  "

$ws.Range("E10").Value = "  -0.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0594"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D14").Value = "1.586.61"
$ws.Range("E14").Value = "  +0.99%  "

$ws.Range("E15").Value = "  +0.74%  "

$ws.Range("E16").Value = "  -0.91%  "

$ws.Range("D17").Value = "28.548.55"
$ws.Range("E17").Value = "  +2.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "0.0₃0695"
$ws.Range("E21").Value = "  -1.11%  "

$ws.Range("E22").Value = "  +0.71%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.60%  "

$ws.Range("E24").Value = "  -0.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("E28").Value = "  -1.31%  "

$ws.Range("E29").Value = "  -1.54%  "

$ws.Range("E30").Value = "  +0.73%  "

$ws.Range("E31").Value = "  -1.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0466"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.35%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("E34").Value = "  +0.24%  "

$ws.Range("D35").Value = "1.397.12"
$ws.Range("E35").Value = "  -0.88%  "

$ws.Range("E36").Value = "  -0.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.05%  "

$ws.Range("E38").Value = "  +2.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.54%  "

$ws.Range("E40").Value = "  -0.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.535"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.58%  "

$ws.Range("E42").Value = "  +0.71%  "

$ws.Range("E43").Value = "  -0.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.979"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("D48").Value = "1.713.92"
$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.85%  "

# Rows 50 and 51 swap: Cronos and BabyDogeCoin exchange positions, with updated values
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0103"
$ws.Range("E50").Value = "  +0.48%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0518"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.35%  "
